$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run that contains document position $pos into two runs
# (a boundary right at $pos) without leaving any formatting residue behind.
# We do this by inserting a paragraph mark at $pos and immediately deleting
# that same mark again - the two resulting runs stay split from then on,
# but (unlike toggling Bold/Italic/etc.) no stray <w:rPr/> is introduced.
# ---------------------------------------------------------------------------
function SplitAt($pos) {
    $r0 = $d.Range($pos, $pos)
    $r0.InsertParagraphAfter()
    $mark = $d.Range($pos, $pos + 1)
    $mark.Delete()
}

$Q1 = [char]0x201C   # “
$Q2 = [char]0x201D   # ”

# ===========================================================================
# Change 1: split "Bobbie-Jo.Webb-Robertson@pnnl.gov" into "bj" + "@pnnl.gov"
# (two runs, same rPr: sz=24 / szCs=24), keeping the neighbouring ", " and
# ")" runs separate, exactly like before the edit.
# ===========================================================================
$rngEmail = $d.Content
$rngEmail.Find.Execute("Bobbie-Jo.Webb-Robertson@pnnl.gov") | Out-Null
$emailStart = $rngEmail.Start
$emailEnd = $rngEmail.End

$fullEmail = $d.Range($emailStart, $emailEnd)
$fullEmail.Text = "bj@pnnl.gov"
$emailNewEnd = $emailStart + ("bj@pnnl.gov".Length)

SplitAt($emailStart)
SplitAt($emailStart + 2)
SplitAt($emailNewEnd)

# ===========================================================================
# Change 2: after "...placed in the "input_seq.fasta" file." (end of the
# "Include new AVP sequences" bullet), append a new sentence about the
# merged training/test file, as three separate runs.
# ===========================================================================
$pInclude = $d.Paragraphs.Item(44)
$insertPos2 = $pInclude.Range.End - 1

$r = $d.Range($insertPos2, $insertPos2)
$r.InsertAfter(" Please note that the training sequences are in the " + $Q1)
$insertPos2 = $insertPos2 + (" Please note that the training sequences are in the " + $Q1).Length

$r = $d.Range($insertPos2, $insertPos2)
$r.InsertAfter("selected_train_test_merged_file.csv")
$insertPos2 = $insertPos2 + ("selected_train_test_merged_file.csv".Length)

$r = $d.Range($insertPos2, $insertPos2)
$r.InsertAfter($Q2 + " file.")

# ===========================================================================
# Change 3: same sentence appended at the end of the "Include new non-AVP
# sequences" bullet (after its trailing " " run), again as three runs.
# ===========================================================================
$pIncludeNon = $d.Paragraphs.Item(45)
$insertPos3 = $pIncludeNon.Range.End - 1

$r = $d.Range($insertPos3, $insertPos3)
$r.InsertAfter("Please note that the training sequences are in the " + $Q1)
$insertPos3 = $insertPos3 + ("Please note that the training sequences are in the " + $Q1).Length

$r = $d.Range($insertPos3, $insertPos3)
$r.InsertAfter("selected_train_test_merged_file.csv")
$insertPos3 = $insertPos3 + ("selected_train_test_merged_file.csv".Length)

$r = $d.Range($insertPos3, $insertPos3)
$r.InsertAfter($Q2 + " file.")

# ===========================================================================
# Change 4: "Citations:" -> "Citation:"
# ===========================================================================
$d.Content.Find.Execute("Citations:", $true, $false, $false, $false, $false, $true, 1, $false, "Citation:", 2) | Out-Null

# ===========================================================================
# Change 5: split "ubmitted. " into "ubmitted." + " ", and rename the paper
# title, keeping the surrounding ". " / "Scientific Reports" / "." runs
# separate the way they were originally.
# ===========================================================================
$oldTitle = "Improved Antiviral Peptide Prediction via Effective Feature Selection"
$newTitle = "Feature-Informed Reduced Machine Learning for Antiviral Peptide Prediction"
$d.Content.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2) | Out-Null

$rngSub = $d.Content
$rngSub.Find.Execute("ubmitted. ") | Out-Null
$subStart = $rngSub.Start

SplitAt($subStart + 9)            # "ubmitted." | " "
SplitAt($subStart + 10)           # " " | new title

$titleStart = $subStart + 10
$titleEnd = $titleStart + $newTitle.Length
SplitAt($titleEnd)                # new title | ". "
SplitAt($titleEnd + 2)            # ". " | "Scientific Reports"
SplitAt($titleEnd + 2 + "Scientific Reports".Length)  # "Scientific Reports" | "."
